$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.039440053278146
$ws.Cells.Item(2, 4).Value = 1.049213022402575
$ws.Cells.Item(2, 5).Value = 1.047921438693354
$ws.Cells.Item(2, 6).Value = 1.058768178272876
$ws.Cells.Item(2, 9).Value = 1.03895710367115
$ws.Cells.Item(2, 10).Value = 1.044532059003349
$ws.Cells.Item(2, 11).Value = 1.051970707526802
$ws.Cells.Item(2, 12).Value = 1.050682725919841
$ws.Cells.Item(2, 13).Value = 1.061499510705602
$ws.Cells.Item(2, 14).Value = 1.018617053227244
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040816588124179
$ws.Cells.Item(3, 4).Value = 1.049898121860066
$ws.Cells.Item(3, 5).Value = 1.049058933976195
$ws.Cells.Item(3, 6).Value = 1.059785210839533
$ws.Cells.Item(3, 9).Value = 1.039106450103648
$ws.Cells.Item(3, 10).Value = 1.045551664742903
$ws.Cells.Item(3, 11).Value = 1.052468336406603
$ws.Cells.Item(3, 12).Value = 1.051631319460515
$ws.Cells.Item(3, 13).Value = 1.062330127060744
$ws.Cells.Item(3, 14).Value = 1.018967428490599
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041707084886664
$ws.Cells.Item(4, 4).Value = 1.050340841974222
$ws.Cells.Item(4, 5).Value = 1.049794908637081
$ws.Cells.Item(4, 6).Value = 1.060442929921757
$ws.Cells.Item(4, 9).Value = 1.039201533388568
$ws.Cells.Item(4, 10).Value = 1.046210772455185
$ws.Cells.Item(4, 11).Value = 1.052789073047302
$ws.Cells.Item(4, 12).Value = 1.052244480799462
$ws.Cells.Item(4, 13).Value = 1.062866612938902
$ws.Cells.Item(4, 14).Value = 1.019193618369989
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04208140236024
$ws.Cells.Item(5, 4).Value = 1.05052682109596
$ws.Cells.Item(5, 5).Value = 1.050104299695271
$ws.Cells.Item(5, 6).Value = 1.060719347622928
$ws.Cells.Item(5, 9).Value = 1.039241134320322
$ws.Cells.Item(5, 10).Value = 1.046487709091636
$ws.Cells.Item(5, 11).Value = 1.052923608458802
$ws.Cells.Item(5, 12).Value = 1.052502101837761
$ws.Cells.Item(5, 13).Value = 1.063091918286945
$ws.Cells.Item(5, 14).Value = 1.0192885832204
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042144249241233
$ws.Cells.Item(6, 4).Value = 1.050558039564441
$ws.Cells.Item(6, 5).Value = 1.050156247125708
$ws.Cells.Item(6, 6).Value = 1.060765754273115
$ws.Cells.Item(6, 9).Value = 1.039247761669382
$ws.Cells.Item(6, 10).Value = 1.046534199109289
$ws.Cells.Item(6, 11).Value = 1.052946179820196
$ws.Cells.Item(6, 12).Value = 1.052545348681302
$ws.Cells.Item(6, 13).Value = 1.063129734391362
$ws.Cells.Item(6, 14).Value = 1.019304520902959
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041712086716564
$ws.Cells.Item(7, 4).Value = 1.050343327588556
$ws.Cells.Item(7, 5).Value = 1.049799042783074
$ws.Cells.Item(7, 6).Value = 1.060446623770892
$ws.Cells.Item(7, 9).Value = 1.039202064000301
$ws.Cells.Item(7, 10).Value = 1.046214473491494
$ws.Cells.Item(7, 11).Value = 1.052790871903834
$ws.Cells.Item(7, 12).Value = 1.05224792373945
$ws.Cells.Item(7, 13).Value = 1.06286962439563
$ws.Cells.Item(7, 14).Value = 1.019194887786821
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.039905306077563
$ws.Cells.Item(8, 4).Value = 1.049444675900495
$ws.Cells.Item(8, 5).Value = 1.048305874281628
$ws.Cells.Item(8, 6).Value = 1.059111966160234
$ws.Cells.Item(8, 9).Value = 1.039007897678967
$ws.Cells.Item(8, 10).Value = 1.04487677452874
$ws.Cells.Item(8, 11).Value = 1.052139144938816
$ws.Cells.Item(8, 12).Value = 1.05100344154766
$ws.Cells.Item(8, 13).Value = 1.0617804245455
$ws.Cells.Item(8, 14).Value = 1.018735573553168
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.036719727842057
$ws.Cells.Item(9, 4).Value = 1.047856670758421
$ws.Cells.Item(9, 5).Value = 1.045674169962287
$ws.Cells.Item(9, 6).Value = 1.056757260081427
$ws.Cells.Item(9, 9).Value = 1.038653854751098
$ws.Cells.Item(9, 10).Value = 1.042514533682928
$ws.Cells.Item(9, 11).Value = 1.050981049546742
$ws.Cells.Item(9, 12).Value = 1.04880550897977
$ws.Cells.Item(9, 13).Value = 1.059853579310382
$ws.Cells.Item(9, 14).Value = 1.017922143425586
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034594565085995
$ws.Cells.Item(10, 4).Value = 1.046795014924574
$ws.Cells.Item(10, 5).Value = 1.043919209398587
$ws.Cells.Item(10, 6).Value = 1.055185464147921
$ws.Cells.Item(10, 9).Value = 1.038409830686645
$ws.Cells.Item(10, 10).Value = 1.040936171021166
$ws.Cells.Item(10, 11).Value = 1.050202482429113
$ws.Cells.Item(10, 12).Value = 1.047336750150355
$ws.Cells.Item(10, 13).Value = 1.058563887235043
$ws.Cells.Item(10, 14).Value = 1.017377086138697
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033673951145631
$ws.Cells.Item(11, 4).Value = 1.046334600718068
$ws.Cells.Item(11, 5).Value = 1.043159147980558
$ws.Cells.Item(11, 6).Value = 1.054504371680625
$ws.Cells.Item(11, 9).Value = 1.03830226946072
$ws.Cells.Item(11, 10).Value = 1.040251851376027
$ws.Cells.Item(11, 11).Value = 1.049863810305201
$ws.Cells.Item(11, 12).Value = 1.046699912717673
$ws.Cells.Item(11, 13).Value = 1.058004205247658
$ws.Cells.Item(11, 14).Value = 1.017140403746029
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.033331928773474
$ws.Cells.Item(12, 4).Value = 1.046163475788541
$ws.Cells.Item(12, 5).Value = 1.042876802028369
$ws.Cells.Item(12, 6).Value = 1.054251307916102
$ws.Cells.Item(12, 9).Value = 1.038262031359642
$ws.Cells.Item(12, 10).Value = 1.03999752975563
$ws.Cells.Item(12, 11).Value = 1.049737779690183
$ws.Cells.Item(12, 12).Value = 1.046463232173303
$ws.Cells.Item(12, 13).Value = 1.057796127249916
$ws.Cells.Item(12, 14).Value = 1.017052388003413
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.033405296765251
$ws.Cells.Item(13, 4).Value = 1.046200187488753
$ws.Cells.Item(13, 5).Value = 1.04293736737833
$ws.Cells.Item(13, 6).Value = 1.054305594400839
$ws.Cells.Item(13, 9).Value = 1.038270675469737
$ws.Cells.Item(13, 10).Value = 1.040052088754706
$ws.Cells.Item(13, 11).Value = 1.049764824219212
$ws.Cells.Item(13, 12).Value = 1.04651400690782
$ws.Cells.Item(13, 13).Value = 1.057840769179549
$ws.Cells.Item(13, 14).Value = 1.01707127227969
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033645680829737
$ws.Cells.Item(14, 4).Value = 1.04632045765189
$ws.Cells.Item(14, 5).Value = 1.043135809706549
$ws.Cells.Item(14, 6).Value = 1.054483454903294
$ws.Cells.Item(14, 9).Value = 1.038298949182373
$ws.Cells.Item(14, 10).Value = 1.040230831838931
$ws.Cells.Item(14, 11).Value = 1.049853397323778
$ws.Cells.Item(14, 12).Value = 1.046680351310915
$ws.Cells.Item(14, 13).Value = 1.057987009282001
$ws.Cells.Item(14, 14).Value = 1.017133130410042
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033793780586217
$ws.Cells.Item(15, 4).Value = 1.046394545947249
$ws.Cells.Item(15, 5).Value = 1.043258073147618
$ws.Cells.Item(15, 6).Value = 1.054593030556071
$ws.Cells.Item(15, 9).Value = 1.038316331773722
$ws.Cells.Item(15, 10).Value = 1.040340943398632
$ws.Cells.Item(15, 11).Value = 1.049907939293689
$ws.Cells.Item(15, 12).Value = 1.04678282419124
$ws.Cells.Item(15, 13).Value = 1.058077087790238
$ws.Cells.Item(15, 14).Value = 1.01717122978838
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.034655654287704
$ws.Cells.Item(16, 4).Value = 1.046825556145166
$ws.Cells.Item(16, 5).Value = 1.043969648759969
$ws.Cells.Item(16, 6).Value = 1.05523065547416
$ws.Cells.Item(16, 9).Value = 1.038416929195867
$ws.Cells.Item(16, 10).Value = 1.040981568312775
$ws.Cells.Item(16, 11).Value = 1.050224926361073
$ws.Cells.Item(16, 12).Value = 1.047378996741012
$ws.Cells.Item(16, 13).Value = 1.058601005310359
$ws.Cells.Item(16, 14).Value = 1.01739277980184
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.035196173365208
$ws.Cells.Item(17, 4).Value = 1.047095727388846
$ws.Cells.Item(17, 5).Value = 1.044415958925608
$ws.Cells.Item(17, 6).Value = 1.055630487367436
$ws.Cells.Item(17, 9).Value = 1.038479523340299
$ws.Cells.Item(17, 10).Value = 1.041383178237075
$ws.Cells.Item(17, 11).Value = 1.050423349527485
$ws.Cells.Item(17, 12).Value = 1.047752729466465
$ws.Cells.Item(17, 13).Value = 1.058929312835287
$ws.Cells.Item(17, 14).Value = 1.017531572561692
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.035511410287925
$ws.Cells.Item(18, 4).Value = 1.047253245172394
$ws.Cells.Item(18, 5).Value = 1.044676269605907
$ws.Cells.Item(18, 6).Value = 1.055863654935512
$ws.Cells.Item(18, 9).Value = 1.038515850397259
$ws.Cells.Item(18, 10).Value = 1.041617345975079
$ws.Cells.Item(18, 11).Value = 1.050538937116076
$ws.Cells.Item(18, 12).Value = 1.047970639087003
$ws.Cells.Item(18, 13).Value = 1.059120689820316
$ws.Cells.Item(18, 14).Value = 1.017612463540877
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.035618891504729
$ws.Cells.Item(19, 4).Value = 1.047306943037791
$ws.Cells.Item(19, 5).Value = 1.044765026437667
$ws.Cells.Item(19, 6).Value = 1.055943151017022
$ws.Cells.Item(19, 9).Value = 1.038528205937695
$ws.Cells.Item(19, 10).Value = 1.041697176870545
$ws.Cells.Item(19, 11).Value = 1.050578324181351
$ws.Cells.Item(19, 12).Value = 1.048044926736088
$ws.Cells.Item(19, 13).Value = 1.059185924242347
$ws.Cells.Item(19, 14).Value = 1.017640034361481
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.035138184822376
$ws.Cells.Item(20, 4).Value = 1.047066747659667
$ws.Cells.Item(20, 5).Value = 1.044368075586305
$ws.Cells.Item(20, 6).Value = 1.055587594125028
$ws.Cells.Item(20, 9).Value = 1.038472826515104
$ws.Cells.Item(20, 10).Value = 1.041340098059023
$ws.Cells.Item(20, 11).Value = 1.050402076043249
$ws.Cells.Item(20, 12).Value = 1.047712640017932
$ws.Cells.Item(20, 13).Value = 1.058894100889256
$ws.Cells.Item(20, 14).Value = 1.017516688089108
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033574895555571
$ws.Cells.Item(21, 4).Value = 1.046285044028175
$ws.Cells.Item(21, 5).Value = 1.043077374105285
$ws.Cells.Item(21, 6).Value = 1.054431081519359
$ws.Cells.Item(21, 9).Value = 1.038290631148447
$ws.Cells.Item(21, 10).Value = 1.040178200201099
$ws.Cells.Item(21, 11).Value = 1.049827321177872
$ws.Cells.Item(21, 12).Value = 1.046631370658457
$ws.Cells.Item(21, 13).Value = 1.057943950395136
$ws.Cells.Item(21, 14).Value = 1.017114917537456
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032591613078863
$ws.Cells.Item(22, 4).Value = 1.045792939259035
$ws.Cells.Item(22, 5).Value = 1.042265711244007
$ws.Cells.Item(22, 6).Value = 1.053703497478201
$ws.Cells.Item(22, 9).Value = 1.0381744283253
$ws.Cells.Item(22, 10).Value = 1.039446887346674
$ws.Cells.Item(22, 11).Value = 1.049464603481981
$ws.Cells.Item(22, 12).Value = 1.045950775979996
$ws.Cells.Item(22, 13).Value = 1.057345469370048
$ws.Cells.Item(22, 14).Value = 1.016861721559137
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.03311290732756
$ws.Cells.Item(23, 4).Value = 1.046053871727283
$ws.Cells.Item(23, 5).Value = 1.042696003905353
$ws.Cells.Item(23, 6).Value = 1.054089245598461
$ws.Cells.Item(23, 9).Value = 1.038236185995438
$ws.Cells.Item(23, 10).Value = 1.039834645179423
$ws.Cells.Item(23, 11).Value = 1.04965701470549
$ws.Cells.Item(23, 12).Value = 1.046311644685197
$ws.Cells.Item(23, 13).Value = 1.057662838667359
$ws.Cells.Item(23, 14).Value = 1.016996001454556
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.03516438746819
$ws.Cells.Item(24, 4).Value = 1.047079842563806
$ws.Cells.Item(24, 5).Value = 1.044389712050972
$ws.Cells.Item(24, 6).Value = 1.055606975882149
$ws.Cells.Item(24, 9).Value = 1.038475853087801
$ws.Cells.Item(24, 10).Value = 1.041359564398234
$ws.Cells.Item(24, 11).Value = 1.050411689076457
$ws.Cells.Item(24, 12).Value = 1.047730754968525
$ws.Cells.Item(24, 13).Value = 1.058910012020471
$ws.Cells.Item(24, 14).Value = 1.017523413940831
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.037543513500947
$ws.Cells.Item(25, 4).Value = 1.048267735607864
$ws.Cells.Item(25, 5).Value = 1.046354606263403
$ws.Cells.Item(25, 6).Value = 1.057366354465459
$ws.Cells.Item(25, 9).Value = 1.038746792826356
$ws.Cells.Item(25, 10).Value = 1.043125842126547
$ws.Cells.Item(25, 11).Value = 1.051281591038233
$ws.Cells.Item(25, 12).Value = 1.049374330718806
$ws.Cells.Item(25, 13).Value = 1.060352614366965
$ws.Cells.Item(25, 14).Value = 1.018132919917558
